{"js": "const body = context.document.body;\n\nconst replacements = [\n    [\"2024-01-29 Monday\", \"2024-01-30 Tuesday\"],\n    [\"11\u00f74=\", \"74\u00f78=\"],\n    [\"12\u00f78=\", \"70\u00f74=\"],\n    [\"96\u00f77=\", \"71\u00f78=\"],\n    [\"85\u00f74=\", \"97\u00f76=\"],\n    [\"42\u00f73=\", \"95\u00f78=\"],\n    [\"98\u00f78=\", \"78\u00f77=\"],\n    [\"69\u00f73=\", \"41\u00f75=\"],\n    [\"37\u00f74=\", \"92\u00f78=\"],\n    [\"42\u00f77=\", \"97\u00f73=\"],\n    [\"82\u00f76=\", \"39\u00f73=\"],\n    [\"91\u00f79=\", \"56\u00f77=\"],\n    [\"82\u00f75=\", \"28\u00f75=\"],\n    [\"38\u00f72=\", \"38\u00f77=\"],\n    [\"97\u00f78=\", \"18\u00f76=\"],\n    [\"53\u00f72=\", \"85\u00f78=\"],\n    [\"55\u00f72=\", \"38\u00f77=\"],\n    [\"46\u00f72=\", \"81\u00f78=\"],\n    [\"67\u00f73=\", \"54\u00f79=\"],\n    [\"27\u00f74=\", \"15\u00f74=\"],\n    [\"21\u00f78=\", \"72\u00f73=\"],\n    [\"31\u00f74=\", \"94\u00f76=\"],\n    [\"73\u00f73=\", \"75\u00f79=\"],\n    [\"96\u00f72=\", \"54\u00f77=\"],\n    [\"71\u00f79=\", \"12\u00f73=\"],\n    [\"88\u00f79=\", \"11\u00f79=\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n    const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length > 0) {\n        results.items[0].insertText(replaceText, \"Replace\");\n        await context.sync();\n    }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-29 Monday\", \"2024-01-30 Tuesday\"),\n    @(\"11\u00f74=\", \"74\u00f78=\"),\n    @(\"12\u00f78=\", \"70\u00f74=\"),\n    @(\"96\u00f77=\", \"71\u00f78=\"),\n    @(\"85\u00f74=\", \"97\u00f76=\"),\n    @(\"42\u00f73=\", \"95\u00f78=\"),\n    @(\"98\u00f78=\", \"78\u00f77=\"),\n    @(\"69\u00f73=\", \"41\u00f75=\"),\n    @(\"37\u00f74=\", \"92\u00f78=\"),\n    @(\"42\u00f77=\", \"97\u00f73=\"),\n    @(\"82\u00f76=\", \"39\u00f73=\"),\n    @(\"91\u00f79=\", \"56\u00f77=\"),\n    @(\"82\u00f75=\", \"28\u00f75=\"),\n    @(\"38\u00f72=\", \"38\u00f77=\"),\n    @(\"97\u00f78=\", \"18\u00f76=\"),\n    @(\"53\u00f72=\", \"85\u00f78=\"),\n    @(\"55\u00f72=\", \"38\u00f77=\"),\n    @(\"46\u00f72=\", \"81\u00f78=\"),\n    @(\"67\u00f73=\", \"54\u00f79=\"),\n    @(\"27\u00f74=\", \"15\u00f74=\"),\n    @(\"21\u00f78=\", \"72\u00f73=\"),\n    @(\"31\u00f74=\", \"94\u00f76=\"),\n    @(\"73\u00f73=\", \"75\u00f79=\"),\n    @(\"96\u00f72=\", \"54\u00f77=\"),\n    @(\"71\u00f79=\", \"12\u00f73=\"),\n    @(\"88\u00f79=\", \"11\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
